$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string need the cell
# format forced to Text first, otherwise Excel auto-converts the entry
# into a numeric value (losing the original text-cell semantics / exact
# decimal representation, e.g. 257.09 -> 257.08999999999997).
$ws.Range('D2').Value = '43.662.33'
$ws.Range('E2').Value = '  +3.26%  '
$ws.Range('D3').Value = '2.187.67'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '257.09'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '80.18'
$ws.Range('E6').Value = '  +8.50%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  +1.76%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.68'
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0917'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.93'
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '2.517.76'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.22'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').Value = '2.170.63'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.772'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '43.680.87'
$ws.Range('E18').Value = '  +3.55%  '
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.05'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.90'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.38'
$ws.Range('E22').Value = '  +11.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '229.75'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('E24').Value = '  -6.21%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '41.07'
$ws.Range('E26').Value = '  +12.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.57'
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.36'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +4.86%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.86'
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.31'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0867'
$ws.Range('E33').Value = '  +7.26%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.24'
$ws.Range('E34').Value = '  +2.91%  '
$ws.Range('E35').Value = '  +7.00%  '
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.42'
$ws.Range('E37').Value = '  +4.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0353'
$ws.Range('E38').Value = '  +5.90%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '13.05'
$ws.Range('E39').Value = '  +11.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.84'
$ws.Range('E40').Value = '  +17.97%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.43'
$ws.Range('E42').Value = '  +5.87%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '61.97'
$ws.Range('E43').Value = '  +4.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.197'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '100.55'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0981'
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.19'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.10'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').Value = '  +3.58%  '
$ws.Range('E50').Value = '  +26.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.437'
$ws.Range('E51').Value = '  -6.85%  '

Write-Output "Applied 103 cell updates"
